$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 89, shifting existing rows 89-130 down to 90-131.
$ws.Rows.Item(89).Insert()

# Match the date-number format used by column D on the other data rows.
$ws.Cells.Item(89, 4).NumberFormat = $ws.Cells.Item(90, 4).NumberFormat

$ws.Cells.Item(89, 1).Value = 10
$ws.Cells.Item(89, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(89, 3).Value = "La Araucanía"
$ws.Cells.Item(89, 4).Value = 44466
$ws.Cells.Item(89, 5).Value = 9
$ws.Cells.Item(89, 6).Value = 100112005
$ws.Cells.Item(89, 7).Value = "Puerro"
$ws.Cells.Item(89, 8).Value = "Azul de Maquehue"
$ws.Cells.Item(89, 9).Value = "Primera"
$ws.Cells.Item(89, 10).Value = 80
$ws.Cells.Item(89, 11).Value = 6000
$ws.Cells.Item(89, 12).Value = 7000
$ws.Cells.Item(89, 13).Value = 6500
$ws.Cells.Item(89, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(89, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(89, 16).Value = 542
$ws.Cells.Item(89, 17).Value = 12
$ws.Cells.Item(89, 18).Value = "Hortaliza"
